$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows 302-328 (dates 2021-06-29 .. 2021-07-25)
$newRows = @(
    @{Row=302; A=44376; B=0; C=0; D=0}
    @{Row=303; A=44377; B=0; C=0; D=0}
    @{Row=304; A=44378; B=0; C=0; D=0}
    @{Row=305; A=44379; B=0; C=0; D=0}
    @{Row=306; A=44380; B=0; C=0; D=0}
    @{Row=307; A=44381; B=0; C=0; D=0}
    @{Row=308; A=44382; B=0; C=0; D=0}
    @{Row=309; A=44383; B=0; C=0; D=0}
    @{Row=310; A=44384; B=0; C=0; D=0}
    @{Row=311; A=44385; B=0; C=0; D=0}
    @{Row=312; A=44386; B=0; C=0; D=0}
    @{Row=313; A=44387; B=0; C=0; D=0}
    @{Row=314; A=44388; B=0; C=0; D=0}
    @{Row=315; A=44389; B=0; C=0; D=0}
    @{Row=316; A=44390; B=0; C=0; D=0}
    @{Row=317; A=44391; B=0; C=0; D=0}
    @{Row=318; A=44392; B=0; C=0; D=0}
    @{Row=319; A=44393; B=0; C=0; D=0}
    @{Row=320; A=44394; B=0; C=0; D=0}
    @{Row=321; A=44395; B=0; C=0; D=0}
    @{Row=322; A=44396; B=0; C=0; D=0}
    @{Row=323; A=44397; B=0; C=0; D=0}
    @{Row=324; A=44398; B=0; C=0; D=0}
    @{Row=325; A=44399; B=1; C=1; D=16.63616702711695}
    @{Row=326; A=44400; B=0; C=1; D=16.63616702711695}
    @{Row=327; A=44401; B=2; C=3; D=49.90850108135086}
    @{Row=328; A=44402; B=0; C=3; D=49.90850108135086}
)

# Use the last existing data row (301) as the style template for column A
# (date-formatted, centered, bordered cells) so no new style entries are created.
$styleTemplate = $ws.Range("A301")

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $styleTemplate.Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
}

$excel.CutCopyMode = 0
